# Automatische test-sync: 2025-06-20 16:30:50
#
# Adds the new "Klacht over levering" mail-log entry (row 19) to the
# "Logs" sheet, bumps its matching category-count row (row 9) on the
# "Dashboard" sheet, and extends the conditional formatting + bar chart
# ranges that were previously capped at the old last row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. "Logs" sheet: append the new incoming-mail row (row 19)
# ---------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A19").Value = "Klacht over levering"
$logs.Range("B19").Value = "mailmind.test@zohomail.eu"
$logs.Range("C19").Value = "Ik ben niet tevreden over mijn bestelling. Ik hoor graag hoe jullie dit oplossen."
$logs.Range("D19").Value = "Klacht / Probleem"
$logs.Range("F19").Value = "2025-06-20 16:30:39"
$logs.Range("G19").Value = "Nee"

# Conditional formatting ranges must grow from row 18 to row 19 as well.
$logs.Range("D2:D18").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D19"))
$logs.Range("G2:G18").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G19"))

# ---------------------------------------------------------------
# 2. "Dashboard" sheet: bump the "Klacht / Probleem" category count
# ---------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A9").Value = "Klacht / Probleem"
$dash.Range("B9").Value = 1

# ---------------------------------------------------------------
# 3. Chart on the Dashboard sheet: extend the series ranges to row 9
# ---------------------------------------------------------------
$chart = $dash.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$9,'Dashboard'!`$B`$2:`$B`$9,1)"
